# Trade #30 closed at 2026-02-17 12:38:15 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refresh aggregate stats after closing trade #30
#  - Strategy Status sheet: refresh MarketMaking strategy row
#  - All Trades / MarketMaking sheets: append the new trade record (row 31)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.72   # Current Capital
$summary.Range("B4").Value = 0.71      # Total P&L $
$summary.Range("B5").Value = 0.47      # Total P&L %
$summary.Range("B6").Value = 30        # Total Trades
$summary.Range("B8").Value = 10        # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.72     # Capital
$status.Range("D4").Value = 30         # Trades
$status.Range("E4").Value = 0.71       # P&L $
$status.Range("F4").Value = 0.72       # P&L %
$status.Range("G4").Value = 40         # Win Rate %

# ---------------------------------------------------------------------------
# 3) Append the new trade row (#30 -> row 31) to both "All Trades" and
#    "MarketMaking" sheets - they carry identical trade logs.
# ---------------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Range("A31").Value = 30
    # Force the date-like text to stay as literal text instead of being
    # auto-converted into a date serial value.
    $ws.Range("B31").NumberFormat = "@"
    $ws.Range("B31").Value = "2026-02-17"
    $ws.Range("C31").Value = "12:38:09"
    $ws.Range("D31").Value = "MarketMaking"
    $ws.Range("E31").Value = "DOWN"
    $ws.Range("F31").Value = 0.25
    $ws.Range("G31").Value = 0.22
    $ws.Range("H31").Value = "CLOSED"
    $ws.Range("I31").Value = -12
    $ws.Range("J31").Value = -0.03
    $ws.Range("K31").Value = 100.72
    $ws.Range("L31").Value = 0
    $ws.Range("M31").Value = 0
    $ws.Range("N31").Value = 0.6
    $ws.Range("O31").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P31").Value = "early_exit"
    $ws.Range("Q31").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking
